$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.8
$ws.Range("H2").Value = 1.74
$ws.Range("I2").Value = 1.77
$ws.Range("K2").Value = 3.8
$ws.Range("N2").Value = 3
$ws.Range("P2").Value = 1.67
$ws.Range("Q2").Value = 2.28
$ws.Range("R2").Value = 1.24
$ws.Range("U2").Value = 1.74
$ws.Range("V2").Value = 2.2
$ws.Range("W2").Value = 1.17
$ws.Range("Z2").Value = 9.199999999999999
$ws.Range("AA2").Value = 18.5
$ws.Range("AB2").Value = 18
$ws.Range("AC2").Value = 8.6
$ws.Range("AH2").Value = 27
$ws.Range("AO2").Value = 16
$ws.Range("F3").Value = 7.8
$ws.Range("G3").Value = 9
$ws.Range("I3").Value = 1.57
$ws.Range("K3").Value = 4.4
$ws.Range("P3").Value = 1.75
$ws.Range("Q3").Value = 2.16
$ws.Range("T3").Value = 2.24
$ws.Range("V3").Value = 2.66
$ws.Range("W3").Value = 1.12
$ws.Range("Y3").Value = 7.4
$ws.Range("Z3").Value = 9.4
$ws.Range("AD3").Value = 12
$ws.Range("AE3").Value = 24
$ws.Range("AJ3").Value = 390
$ws.Range("G4").Value = 3.85
$ws.Range("I4").Value = 2.48
$ws.Range("M4").Value = 1.09
$ws.Range("N4").Value = 3.15
$ws.Range("S4").Value = 3.95
$ws.Range("T4").Value = 1.85
$ws.Range("U4").Value = 1.97
$ws.Range("V4").Value = 1.67
$ws.Range("W4").Value = 1.35
$ws.Range("X4").Value = 12.5
$ws.Range("Y4").Value = 9.199999999999999
$ws.Range("Z4").Value = 15
$ws.Range("AA4").Value = 34
$ws.Range("AB4").Value = 12.5
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 12
$ws.Range("AE4").Value = 29
$ws.Range("AF4").Value = 26
$ws.Range("AG4").Value = 15.5
$ws.Range("AH4").Value = 20
$ws.Range("AI4").Value = 46
$ws.Range("AN4").Value = 55
$ws.Range("F5").Value = 1.64
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 7.2
$ws.Range("V5").Value = 1.16
$ws.Range("W5").Value = 2.42
$ws.Range("AD5").Value = 27
$ws.Range("AL5").Value = 48
$ws.Range("AM5").Value = 200
$ws.Range("AN5").Value = 13
$ws.Range("AO5").Value = 190
$ws.Range("F6").Value = 2.3
$ws.Range("G6").Value = 2.44
$ws.Range("H6").Value = 3.8
$ws.Range("I6").Value = 4.2
$ws.Range("J6").Value = 2.98
$ws.Range("M6").Value = 1.14
$ws.Range("N6").Value = 2.4
$ws.Range("O6").Value = 1.62
$ws.Range("P6").Value = 1.47
$ws.Range("Q6").Value = 2.88
$ws.Range("S6").Value = 6.2
$ws.Range("T6").Value = 2.26
$ws.Range("V6").Value = 1.32
$ws.Range("W6").Value = 1.69
$ws.Range("Y6").Value = 10
$ws.Range("Z6").Value = 980
$ws.Range("AB6").Value = 7.8
$ws.Range("AC6").Value = 8.199999999999999
$ws.Range("AH6").Value = 980
$ws.Range("AJ6").Value = 980
$ws.Range("AK6").Value = 980
$ws.Range("AM6").Value = 290
$ws.Range("G7").Value = 1.99
$ws.Range("H7").Value = 5.1
$ws.Range("I7").Value = 5.9
$ws.Range("J7").Value = 3.2
$ws.Range("M7").Value = 1.13
$ws.Range("P7").Value = 1.53
$ws.Range("Q7").Value = 2.62
$ws.Range("R7").Value = 1.19
$ws.Range("S7").Value = 5.4
$ws.Range("U7").Value = 1.65
$ws.Range("V7").Value = 1.21
$ws.Range("X7").Value = 8.6
$ws.Range("Y7").Value = 14
$ws.Range("AB7").Value = 6.6
$ws.Range("AC7").Value = 8
$ws.Range("AF7").Value = 10
$ws.Range("AG7").Value = 11.5
$ws.Range("I8").Value = 3.8
$ws.Range("J8").Value = 3.1
$ws.Range("O8").Value = 1.48
$ws.Range("F9").Value = 1.71
$ws.Range("G9").Value = 1.78
$ws.Range("H9").Value = 5.5
$ws.Range("I9").Value = 6.2
$ws.Range("J9").Value = 3.8
$ws.Range("M9").Value = 1.08
$ws.Range("P9").Value = 1.8
$ws.Range("R9").Value = 1.3
$ws.Range("S9").Value = 3.75
$ws.Range("T9").Value = 1.99
$ws.Range("U9").Value = 1.84
$ws.Range("V9").Value = 1.19
$ws.Range("W9").Value = 2.26
$ws.Range("X9").Value = 13.5
$ws.Range("Y9").Value = 18
$ws.Range("AB9").Value = 7.8
$ws.Range("AD9").Value = 24
$ws.Range("AF9").Value = 9.800000000000001
$ws.Range("AG9").Value = 10.5
$ws.Range("AJ9").Value = 18.5
$ws.Range("AL9").Value = 980
$ws.Range("AN9").Value = 13
$ws.Range("F10").Value = 1.99
$ws.Range("J10").Value = 3.2
$ws.Range("K10").Value = 3.4
$ws.Range("M10").Value = 1.12
$ws.Range("Q10").Value = 2.56
$ws.Range("S10").Value = 5.2
$ws.Range("U10").Value = 1.74
$ws.Range("V10").Value = 1.26
$ws.Range("W10").Value = 1.9
$ws.Range("X10").Value = 8.800000000000001
$ws.Range("AE10").Value = 110
$ws.Range("F11").Value = 1.59
$ws.Range("G11").Value = 1.62
$ws.Range("J11").Value = 4.5
$ws.Range("K11").Value = 4.7
$ws.Range("L11").Value = 1.32
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 4.9
$ws.Range("O11").Value = 1.24
$ws.Range("P11").Value = 2.32
$ws.Range("Q11").Value = 1.72
$ws.Range("R11").Value = 1.53
$ws.Range("S11").Value = 2.72
$ws.Range("T11").Value = 1.8
$ws.Range("U11").Value = 2.16
$ws.Range("W11").Value = 2.6
$ws.Range("X11").Value = 23
$ws.Range("Y11").Value = 23
$ws.Range("AB11").Value = 10
$ws.Range("AC11").Value = 10
$ws.Range("AE11").Value = 85
$ws.Range("AF11").Value = 10.5
$ws.Range("AH11").Value = 21
$ws.Range("AI11").Value = 990
$ws.Range("AJ11").Value = 15.5
$ws.Range("AK11").Value = 15.5
$ws.Range("AM11").Value = 1000
$ws.Range("AO11").Value = 80
